$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$ax2 = $chart.Axes(2)
$title = $ax2.AxisTitle
Write-Host "Orientation before:" $title.Orientation
$title.Text = "Average Accuracy"
Write-Host "Orientation after:" $title.Orientation
